$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New K (Strikeout) values computed from regenerated save_data (column G, rows 2-16)
$kValues = @{
    2  = 4
    3  = 2
    4  = 2
    5  = 2
    6  = 2
    7  = 2
    8  = 1
    9  = 2
    10 = 0
    11 = 1
    12 = 2
    13 = 0
    14 = 3
    15 = 2
    16 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
